$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 10.92359866666667
$ws.Range("N2").Value = 32.770796
$ws.Range("O2").Value = 0.2236009040380497
$ws.Range("P2").Value = 0.2236009040380497
$ws.Range("Q2").Value = 0.800517722288889
$ws.Range("R2").Value = 7.204659500600001
$ws.Range("S2").Value = 0.2236009040380497
$ws.Range("T2").Value = 0.2236009040380497

# Row 3
$ws.Range("O3").Value = 0.4261214970992155
$ws.Range("P3").Value = 0.4261214970992155
$ws.Range("S3").Value = 0.4261214970992155
$ws.Range("T3").Value = 0.4261214970992155

# Row 4
$ws.Range("M4").Value = 13.06524766666667
$ws.Range("N4").Value = 39.195743
$ws.Range("O4").Value = 0.2674394472823625
$ws.Range("P4").Value = 0.2674394472823625
$ws.Range("Q4").Value = 0.9574648998388887
$ws.Range("R4").Value = 8.61718409855
$ws.Range("S4").Value = 0.2674394472823625
$ws.Range("T4").Value = 0.2674394472823625

# Row 5
$ws.Range("M5").Value = 4.046901
$ws.Range("N5").Value = 12.140703
$ws.Range("O5").Value = 0.0828381515803724
$ws.Range("P5").Value = 0.0828381515803724
$ws.Range("Q5").Value = 0.29657039495
$ws.Range("R5").Value = 2.66913355455
$ws.Range("S5").Value = 0.0828381515803724
$ws.Range("T5").Value = 0.0828381515803724
